$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value (from meteocat daily-summary refresh).
# Percentage cells (column H) are pre-formatted as Text ("@") before assignment so
# Excel stores the literal "NN%" string instead of auto-converting it to a percent number.
$updates = @(
    @{ Cell = "E2"; Value = "2026-02-21 18:18:15" }
    @{ Cell = "E3"; Value = "2026-02-21 18:18:17" }
    @{ Cell = "E4"; Value = "2026-02-21 18:18:19" }
    @{ Cell = "H4"; Value = "70%" }
    @{ Cell = "E5"; Value = "2026-02-21 18:18:21" }
    @{ Cell = "O5"; Value = "3.5 °C" }
    @{ Cell = "E6"; Value = "2026-02-21 18:18:24" }
    @{ Cell = "E7"; Value = "2026-02-21 18:18:26" }
    @{ Cell = "E8"; Value = "2026-02-21 18:18:28" }
    @{ Cell = "E9"; Value = "2026-02-21 18:18:31" }
    @{ Cell = "H9"; Value = "51%" }
    @{ Cell = "O9"; Value = "14.2 °C" }
    @{ Cell = "E10"; Value = "2026-02-21 18:18:32" }
    @{ Cell = "E11"; Value = "2026-02-21 18:18:33" }
    @{ Cell = "E12"; Value = "2026-02-21 18:18:34" }
    @{ Cell = "H12"; Value = "57%" }
    @{ Cell = "E13"; Value = "2026-02-21 18:18:35" }
    @{ Cell = "H13"; Value = "61%" }
    @{ Cell = "O13"; Value = "5.2 °C" }
    @{ Cell = "E14"; Value = "2026-02-21 18:18:36" }
    @{ Cell = "H14"; Value = "66%" }
    @{ Cell = "E15"; Value = "2026-02-21 18:18:37" }
    @{ Cell = "H15"; Value = "50%" }
    @{ Cell = "O15"; Value = "14.0 °C" }
    @{ Cell = "E16"; Value = "2026-02-21 18:18:39" }
    @{ Cell = "E17"; Value = "2026-02-21 18:18:40" }
    @{ Cell = "E18"; Value = "2026-02-21 18:18:41" }
    @{ Cell = "O18"; Value = "8.8 °C" }
    @{ Cell = "E19"; Value = "2026-02-21 18:18:42" }
    @{ Cell = "O19"; Value = "8.0 °C" }
    @{ Cell = "E20"; Value = "2026-02-21 18:18:44" }
    @{ Cell = "E21"; Value = "2026-02-21 18:18:47" }
    @{ Cell = "H21"; Value = "56%" }
    @{ Cell = "J21"; Value = "1030.6 hPa" }
    @{ Cell = "O21"; Value = "7.5 °C" }
    @{ Cell = "E22"; Value = "2026-02-21 18:18:49" }
    @{ Cell = "H22"; Value = "33%" }
    @{ Cell = "E23"; Value = "2026-02-21 18:18:51" }
    @{ Cell = "O23"; Value = "2.7 °C" }
    @{ Cell = "E24"; Value = "2026-02-21 18:18:54" }
    @{ Cell = "J24"; Value = "1031.6 hPa" }
    @{ Cell = "O24"; Value = "6.5 °C" }
    @{ Cell = "E25"; Value = "2026-02-21 18:18:56" }
    @{ Cell = "E26"; Value = "2026-02-21 18:18:58" }
    @{ Cell = "E27"; Value = "2026-02-21 18:19:01" }
    @{ Cell = "E28"; Value = "2026-02-21 18:19:03" }
    @{ Cell = "O28"; Value = "8.3 °C" }
    @{ Cell = "E29"; Value = "2026-02-21 18:19:05" }
    @{ Cell = "H29"; Value = "63%" }
    @{ Cell = "E30"; Value = "2026-02-21 18:19:08" }
    @{ Cell = "H30"; Value = "63%" }
    @{ Cell = "O30"; Value = "12.0 °C" }
    @{ Cell = "E31"; Value = "2026-02-21 18:19:10" }
    @{ Cell = "K31"; Value = "15.1 MJ/m2" }
    @{ Cell = "E32"; Value = "2026-02-21 18:19:13" }
    @{ Cell = "E33"; Value = "2026-02-21 18:19:15" }
    @{ Cell = "H33"; Value = "54%" }
    @{ Cell = "J33"; Value = "1030.3 hPa" }
    @{ Cell = "O33"; Value = "6.4 °C" }
    @{ Cell = "E34"; Value = "2026-02-21 18:19:18" }
    @{ Cell = "E35"; Value = "2026-02-21 18:19:20" }
    @{ Cell = "E36"; Value = "2026-02-21 18:19:22" }
    @{ Cell = "H36"; Value = "54%" }
    @{ Cell = "K36"; Value = "15.2 MJ/m2" }
    @{ Cell = "E37"; Value = "2026-02-21 18:19:25" }
    @{ Cell = "E38"; Value = "2026-02-21 18:19:27" }
    @{ Cell = "O38"; Value = "9.9 °C" }
    @{ Cell = "E39"; Value = "2026-02-21 18:19:30" }
    @{ Cell = "E40"; Value = "2026-02-21 18:19:32" }
    @{ Cell = "J40"; Value = "1030.3 hPa" }
    @{ Cell = "E41"; Value = "2026-02-21 18:19:34" }
    @{ Cell = "O41"; Value = "11.5 °C" }
    @{ Cell = "E42"; Value = "2026-02-21 18:19:36" }
    @{ Cell = "H42"; Value = "71%" }
    @{ Cell = "E43"; Value = "2026-02-21 18:19:39" }
    @{ Cell = "H43"; Value = "76%" }
    @{ Cell = "O43"; Value = "6.8 °C" }
    @{ Cell = "E44"; Value = "2026-02-21 18:19:41" }
    @{ Cell = "E45"; Value = "2026-02-21 18:19:43" }
    @{ Cell = "E46"; Value = "2026-02-21 18:19:46" }
    @{ Cell = "O46"; Value = "10.1 °C" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Cell -match "^H\d+$") {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
